# "Add files via upload" - re-upload of hospital_death_01.xlsx with refreshed
# blood_pressure (column H) lookups and some cosmetic view/column-width tweaks.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Column H (blood_pressure) values were refreshed ------------------
# Rows 2-58 (minus 27 & 52, which were already #N/A) now resolve to #N/A.
$naRows = @(2,3,4,5,6,7,8,9,10,11,12,13,14,15,16,17,18,19,20,21,22,23,24,25,26,28,29,30,31,32,33,34,35,36,37,38,39,40,41,42,43,44,45,46,47,48,49,50,51,53,54,55,56,57,58)
foreach ($row in $naRows) {
    $ws.Cells.Item($row, 8).Value = "#N/A"
}

# Rows 59-130 now resolve to different (refreshed) blood-pressure readings.
$newVals = @{59=163;60=163;61=149;62=90;63=134;64=110;65=160;66=89;67=89;68=160;69=99;70=90;71=160;72=103;73=105;74=130;75=112;76=112;77=160;78=137;79=137;80=137;81=137;82=137;83=137;84=137;85=137;86=137;87=137;88=137;89=137;90=137;91=137;92=137;93=152;94=160;95=146;96=90;97=121;98=121;99=160;100=160;101=160;102=160;103=154;104=160;105=160;106=160;107=160;108=140;109=177;110=160;111=114;112=106;113=143;114=160;115=123;116=155;117=120;118=104;119=112;120=160;121=160;122=133;123=95;124=110;125=110;126=110;127=102;128=93;129=102;130=115}
foreach ($row in $newVals.Keys) {
    $ws.Cells.Item($row, 8).Value = $newVals[$row]
}

# --- 2. Column width tweaks ------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 25.166666666666668
$ws.Columns.Item(2).ColumnWidth = 8.43
$ws.Columns.Item(3).ColumnWidth = 8.916666666666666
$ws.Columns.Item(4).ColumnWidth = 15.916666666666666
$ws.Columns.Item(5).ColumnWidth = 8.43

# --- 3. Defined name: _FilterDatabase range lost its last column ----------
$wb.Names.Item("_xlnm._FilterDatabase").RefersTo = "=גיליון1!`$A`$1:`$F`$130"

# --- 4. Selection / scroll position moved ----------------------------------
$ws.Range("B9").Select()
